$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 19 (shifts existing rows 19..109 down to 20..110)
$ws.Rows("19:19").Insert()

# Populate the new row 19 with the new weekly record
$ws.Range("A19").Value = 8
$ws.Range("B19").Value = "Terminal La Palmera de La Serena"
$ws.Range("C19").Value = "Coquimbo"
$ws.Range("D19").Value = 44560
$ws.Range("E19").Value = 4
$ws.Range("F19").Value = 100112044
$ws.Range("G19").Value = "Perejil"
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 3400
$ws.Range("K19").Value = 2500
$ws.Range("L19").Value = 3000
$ws.Range("M19").Value = 2750
$ws.Range("N19").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O19").Value = "Provincia del Elquí"
$ws.Range("P19").Value = 1833
$ws.Range("Q19").Value = 1.5
$ws.Range("R19").Value = "Hortaliza"
